# Insert a new localisation row ("fb_logout_ask") above the former row 19
# ("quit_ask") on Sheet1, pushing every row below it down by one - this is
# the "replace background menu, start button menu" commit: it adds a new
# Facebook-logout confirmation string set that needs a table row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a fresh row at position 19; Excel shifts rows 19:45 down to 20:46
# and copies formatting from the row above, which already matches the
# target styles (A=6, B=1, C=1).
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row with the new Id / English / Vietnamese
# strings.
$ws.Range("A19").Value = "fb_logout_ask"
$ws.Range("B19").Value = "Do you want to disconnect Facebook? "
$ws.Range("C19").Value = "Bạn có muốn ngắt kết nối Facebook!"

# Match the new selection / scroll position recorded in the workbook view.
$ws.Range("B18").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Top = 900
